# Updated main GSC export data: append 4 new daily rows (2025-12-26 .. 2025-12-29)
# to the "Chart" sheet, mirroring the existing row layout (Date | No video indexed |
# Video indexed | Impressions). The last two new days don't have impressions data
# yet, so column D is written as an (empty) text value, matching the existing
# "no data yet" convention used elsewhere in this export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

function Set-TextCell($range, [string]$text) {
    # Writing a leading single-quote forces Excel to treat the value as literal
    # text instead of auto-converting date-shaped strings (e.g. "2025-12-26")
    # into date serials. ClearFormats() then drops the transient "quote prefix"
    # style so the cell keeps the sheet's normal (default) style, same as every
    # other cell in this column.
    $range.Value = "'" + $text
    $range.ClearFormats()
}

$rows = @(
    @{ Row = 84; Date = "2025-12-26"; NoVideo = 22; Video = 1; Impressions = 0 },
    @{ Row = 85; Date = "2025-12-27"; NoVideo = 22; Video = 1; Impressions = 0 },
    @{ Row = 86; Date = "2025-12-28"; NoVideo = 22; Video = 1; Impressions = $null },
    @{ Row = 87; Date = "2025-12-29"; NoVideo = 22; Video = 1; Impressions = $null }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    Set-TextCell $ws.Range("A$rowNum") $r.Date
    $ws.Range("B$rowNum").Value = $r.NoVideo
    $ws.Range("C$rowNum").Value = $r.Video
    if ($null -eq $r.Impressions) {
        # No impressions value recorded yet for this day -> empty text cell
        # (same representation as the rest of the export uses for "no data").
        Set-TextCell $ws.Range("D$rowNum") ""
    } else {
        $ws.Range("D$rowNum").Value = $r.Impressions
    }
}
